# Update the "Note" column (H) with the full airline history text that
# replaces the generic "Do you have any extra history..." placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "Ada Air was formed by Albanian/French Ada Group as Adalbanair n 1991 and started operations on 03Feb1992`n" + [char]0x2013 + " Regional carrier operating scheduled services from Albania to holiday destinations"

$ws.Range("H3").Value = "Ada Air was formed by Albanian/French Ada Group as Adalbanair n 1991 and started operations on 03Feb1992"

$ws.Range("H4").Value = "Air Albania was founded on 16/May/2018. The carrier had planned to start operations on 31/Aug/2018 but was delayed. It operated it" + [char]0x2019 + "s first test flight on 15/Sep/2018 using an Airbus A319-100. Received its AOC in Mar2019 and started operations on 19/Apr/2019`n< operates scheduled air services for passengers, based out of Tirana Rinas Mother Teresa (TIA) >"

$ws.Range("H5").Value = "Albanian Airlines (i) was formed as joint venture between Albtransport and Austria" + [char]0x2019 + "s Tyrolean Airways in May1991 under the name of Arberia Airlines. Renamed to current name in May1992. The carrier suspended operations in 1994 and reformed in 1995, starting operations on 20Jun1995. Eventually ceased operations on 10Nov2011 when AOC was revoked`n~ was national airline of Albania operating scheduled services out of Tirana"

$ws.Range("H6").Value = "Albanian Airlines (ii) is a planned scheduled carrier based out of Tirana Rinas Mother Teresa (TIA) & is expecting to start operations in 2017"

$ws.Range("H7").Value = "Albanian Airways was a planned scheduled carrier to be based out of Tirana Rinas Mother Teresa (TIA). Unfortunately it failed to start operations"

$ws.Range("H8").Value = "Albatros Airways was established and started operations on 03Nov2004. On 01Sep2006 the airline was grounded"

$ws.Range("H9").Value = "Albawings was founded in Feb2015 and was awarded its Air Operator" + [char]0x2018 + "s Certificate (AOC) by the Albanian Civil Aviation Authority on 04Feb2016. It started operations in Sep2016`n< operates flights from its hub primarily to several Italian destinations, but also offers flights to destinations in Germany, Slovakia, and the United Kingdom. Low-cost carrier based out of Tirana International Airport >"

$ws.Range("H10").Value = "Albtransport (Altes Transport Tirana) was based in Tirana and eventually became a handling company."

$ws.Range("H11").Value = "Nil further at this time`n" + [char]0x2013 + " [Albanian Airlines (i) was formed as joint venture between Albtransport and Austria" + [char]0x2019 + "s Tyrolean Airways in May1991 under the name of Arberia Airlines. Renamed to current name in May1992. The carrier suspended operations in 1994 and reformed in 1995, starting operations on 20Jun1995. Eventually ceased operations on 10Nov2011 when AOC was revoked`n~ was national airline of Albania operating scheduled services out of Tirana]"

$ws.Range("H12").Value = "Belle Air was founded in 2005 & ceased operations on 24Nov2013 stating economic difficulties (Scheduled low-cost carrier based in Tirana)"

$ws.Range("H13").Value = "No further information at this time"

$ws.Range("H14").Value = "[Ala Littoria was founded in 1923 as AERO EXPRESSO ITALIANA. In 1934, the airline merged with NAVIGAZIONE AEREA " + [char]0x2013 + " SANA, SISA ( Societ" + [char]0x00E0 + " Italiana Servizi Aerei) and SAM (Societ" + [char]0x00E0 + " Aerea Mediterranea), the pre-war national airline of Italy. In 1935, SOCIETA ADRIA AERO LLOYD in Albania (then de facto an Italian colony) and NORD AFRICA AVIAZIONE, operating in the Italian colonies of Libya were integrated in. The airline was used by the Italian military at the start of World War II and ceased operations]"

$ws.Range("H15").Value = "Star Airways was established in Tirana by Italian investors. ItAli Airlines was contracted to operate on Star Airways" + [char]0x2019 + " behalf from 07Jun2010, launching services to Milan Malpensa and Pisa. In Jul2010 the Albanian Civil Aviation Authority suspended Star Airways" + [char]0x2019 + " air operator" + [char]0x2019 + "s certificate. Services continued between Albania and Italy only, using the aircraft and air operator" + [char]0x2019 + "s certificate of ItAli Airlines"

$ws.Range("H16").Value = "Tafa Air was formed in 2009 & commenced operations 18Dec2009. Ceased in early Feb2010 with scheduled flights out of Tirana International Airport and Pristina International Airport to Athens International Airport"
